$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15: Registered Companies (fill E first, then A, then D, matching source order)
$ws.Range("E15").Value = "http://economictimes.indiatimes.com/news/economy/policy/number-of-registered-companies-climbs-to-15-27-lakh-in-january/articleshow/51026331.cms"
$ws.Range("A15").Value = "Registered Companies"
$ws.Range("B15").Value = 15.27
$ws.Range("C15").Value = 100000
$ws.Range("C15").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D15").Value = "companies"

# Row 16: Active Companies
$ws.Range("A16").Value = "Active Companies"
$ws.Range("B16").Value = 10.76
$ws.Range("C16").Value = 100000

# Row 17: petty cash Min
$ws.Range("A17").Value = "petty cash Min"
$ws.Range("B17").Value = 6000
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "Rs"
$ws.Range("E17").Value = "http://www.quickbooks.in/r/accounting-taxes/creating-a-petty-cash-system/"

# Row 18: petty cash Max
$ws.Range("A18").Value = "petty cash Max"
$ws.Range("B18").Value = 30000
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "Rs"

$ws.Range("A17").Select()
